$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.023.11'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '3.765.53'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'630.75"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.18%  '
$ws.Range('D6').Value = "'165.77"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('D7').Value = '3.762.58'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').Value = '  +0.15%  '
$ws.Range('E10').Value = '  -2.17%  '
$ws.Range('D11').Value = "'0.457"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = "'6.76"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.29%  '
$ws.Range('D13').Value = "'0.0000239"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.45%  '
$ws.Range('D14').Value = "'34.86"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('D15').Value = '4.396.76'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '3.762.05'
$ws.Range('E16').Value = '  +0.68%  '
$ws.Range('D17').Value = '68.992.23'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('E18').Value = '  -3.29%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('D20').Value = "'7.02"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').Value = "'461.94"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.29%  '
$ws.Range('D22').Value = "'9.49"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.70%  '
$ws.Range('D23').Value = "'0.704"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = "'0.0000144"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.85%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'82.04"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.23%  '
$ws.Range('E26').Value = '  +0.26%  '
$ws.Range('D27').Value = "'2.12"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.75%  '
$ws.Range('D28').Value = "'10.12"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '3.918.02'
$ws.Range('E30').Value = '  -1.47%  '
$ws.Range('D31').Value = "'2.28"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.93%  '
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('D33').Value = "'7.06"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.55%  '
$ws.Range('D34').Value = "'28.35"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.16%  '
$ws.Range('D35').Value = "'0.174"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +17.47%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '3.716.14'
$ws.Range('E37').Value = '  -1.55%  '
$ws.Range('D38').Value = "'8.89"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = "'3.30"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.65%  '
$ws.Range('D41').Value = "'5.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.14%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = "'0.958"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.38%  '
$ws.Range('D45').Value = "'156.79"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = "'1.97"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = "'47.02"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('D49').Value = "'42.74"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').Value = "'8.34"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.94%  '
